# invitees_test.xlsx update:
#   Row 2 previously held Torgeir's info (email hyperlinked in A2, name in B2).
#   It is replaced with Kjetil's info: email in A2 (plain text, no hyperlink),
#   name "Kjetil" in B2. The obsolete "Hyperlink" cell style is removed and the
#   active selection moves to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto: hyperlink that lived on A2.
$ws.Hyperlinks.Delete() | Out-Null

# A2 currently carries the built-in "Hyperlink" formatting (underline, themed
# color). Re-base its direct formatting on B2's plain style so it again looks
# like an ordinary text cell once the hyperlink is gone.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Swap in the new invitee's data.
$ws.Range("A2").Value = "gotvassli@gmail.com"
$ws.Range("B2").Value = "Kjetil"

# The named "Hyperlink" cell style is no longer referenced anywhere - remove it.
$wb.Styles.Item("Hyperlink").Delete() | Out-Null

# Match the saved selection state (B2 active).
$ws.Range("B2").Select() | Out-Null
